$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "244.72"
Set-TextValue "D4" "5.390"
Set-TextValue "D5" "0.06036"
Set-TextValue "D6" "3.394"
Set-TextValue "D7" "0.8148"
Set-TextValue "D8" "0.9323"
Set-TextValue "D9" "0.1434"
Set-TextValue "D10" "0.07508"
Set-TextValue "D11" "0.03445"
Set-TextValue "D12" "0.03052"
Set-TextValue "D13" "0.09430"
Set-TextValue "D14" "4.016"
Set-TextValue "D15" "0.001590"
Set-TextValue "D16" "0.04811"
Set-TextValue "D17" "0.0005943"
Set-TextValue "D18" "0.005596"
Set-TextValue "D19" "0.004164"
Set-TextValue "D20" "0.0009907"
Set-TextValue "D21" "3.664"
Set-TextValue "D22" "6.451"
Set-TextValue "D23" "2.181"
Set-TextValue "D26" "0.00008403"
Set-TextValue "D40" "0.03994"
Set-TextValue "D41" "0.1078"
Set-TextValue "D43" "0.003052"
Set-TextValue "D44" "0.005783"
Set-TextValue "D45" "0.00005255"
Set-TextValue "D48" "0.002327"
Set-TextValue "D49" "0.00002101"
